$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Pablo Romero"
$ws.Range("B3").Value = 573114524438
$ws.Range("C3").Value = "Caballero"

$ws.Range("C3").Select()
